$wb = $excel.ActiveWorkbook

$wsNotes = $wb.Worksheets.Item("Notes")
$wsData  = $wb.Worksheets.Item("Data")

# ---- Notes sheet: fix the "Units of measure" line ----
$wsNotes.Range("A3").Value = "Units of measure: constant 2015 US$"

# ---- Data sheet: populate the country rows ----
$rows = @(
    @("AG", "Antigua & Barbuda", 2015, 48596.3),
    @("BZ", "Belize", 2015, 4224337),
    @("DM", "Dominica", 2015, 3471354),
    @("GD", "Grenada", 2015, 7267653),
    @("GY", "Guyana", 2015, 6032300),
    @("HT", "Haiti", 2015, 9583210),
    @("JM", "Jamaica", 2015, 12960970),
    @("MS", "Montserrat", 2015, 1904593),
    @("north-central-america", "North & Central America, regional", 2015, 0),
    @("LC", "Saint Lucia", 2015, 8352260),
    @("VC", "Saint Vincent & the Grenadines", 2015, 3301664),
    @("SR", "Suriname", 2015, 0)
)

$r = 2
foreach ($row in $rows) {
    $wsData.Cells.Item($r, 1).Value = $row[0]
    $wsData.Cells.Item($r, 2).Value = $row[1]
    $wsData.Cells.Item($r, 3).Value = $row[2]
    $wsData.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
